$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.563.17"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.73%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.454.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -6.42%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.90%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.36"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.85%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.453.30"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -6.40%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.475"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.16%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.83"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.417"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.90%  "

# Row 13
$ws.Range("E13").Value = "  -7.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.043.56"
$ws.Range("D14").ClearFormats()

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.449.17"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.61%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.534.04"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -4.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.83"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.50%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "438.63"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.99%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.94"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -14.25%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.618"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.07%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.597.20"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.33%  "

# Row 27
$ws.Range("E27").Value = "  -3.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.81%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.18"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -11.25%  "

# Row 30
$ws.Range("E30").Value = "  -6.13%  "

# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.53"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -10.61%  "

# Row 33
$ws.Range("E33").Value = "  -5.43%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.07"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.20%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.444.49"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.81"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.54%  "

# Row 39
$ws.Range("E39").Value = "  +0.15%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.54%  "

# Row 42
$ws.Range("E42").Value = "  -5.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0854"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.66%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.37"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.97%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.874"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.31%  "

# Row 46
$ws.Range("E46").Value = "  -3.74%  "

# Row 47
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.51%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.76"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -12.30%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.49"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.85%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.45"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -13.77%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.994"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.89%  "
